# Daily attendance processing - 2026-01-22 08:06:28
# Re-normalizes the "Recorded By" (column G) lists on the active sheet:
# the automated recorder name ("System"/"system") used to be listed before
# the human recorder's account; it is now moved to the end of the
# comma-separated list of recorders.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $value = $cell.Value2

    if ($value -eq $null) { continue }

    $text = [string]$value

    if ($text -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    }
    elseif ($text -eq "system, backup@backdoor.com, System") {
        $cell.Value = "backup@backdoor.com, System, system"
    }
}
